$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160. This shifts the existing rows 160-201
# down to 161-202 (and extends the used range/dimension accordingly).
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with a new data record (same shape
# as the surrounding rows: Feria Lagunitas de Puerto Montt / Los Lagos /
# Zapallo / Camote, but with its own date, quality grade and origin).
$ws.Cells.Item(160, 1).Value = 4
$ws.Cells.Item(160, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(160, 3).Value = "Los Lagos"
$ws.Cells.Item(160, 4).Value = 44511
$ws.Cells.Item(160, 5).Value = 10
$ws.Cells.Item(160, 6).Value = 100112045
$ws.Cells.Item(160, 7).Value = "Zapallo"
$ws.Cells.Item(160, 8).Value = "Camote"
$ws.Cells.Item(160, 9).Value = "2a nueva(o)"
$ws.Cells.Item(160, 10).Value = 300
$ws.Cells.Item(160, 11).Value = 600
$ws.Cells.Item(160, 12).Value = 600
$ws.Cells.Item(160, 13).Value = 600
$ws.Cells.Item(160, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(160, 15).Value = "Perú"
$ws.Cells.Item(160, 16).Value = 600
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"
